# "fix: delete VNS, perserve best desent"
#
# The "VNS" results (column C, variable-neighbourhood-search best found
# solution) are dropped for the last remaining batch of rows (row 33 and
# rows 242-271) that still carried them - mirroring every other data row,
# which already has no column C value. Column C for rows 2-31 is refreshed
# with the "best descent" run's numbers instead of being deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh rows 2-31 (best-descent column C values) -----------------
$newC = @{
    2  = 23315
    3  = 22748
    4  = 22173
    5  = 23482
    6  = 23160
    7  = 25021
    8  = 24551
    9  = 23443
    10 = 23550
    11 = 23005
    12 = 43035
    13 = 40819
    14 = 40509
    15 = 42453
    16 = 40613
    17 = 39423
    18 = 41206
    19 = 42286
    20 = 41867
    21 = 43605
    22 = 60277
    23 = 63154
    24 = 59483
    25 = 61255
    26 = 61284
    27 = 59326
    28 = 62792
    29 = 61562
    30 = 60302
    31 = 60322
}

foreach ($r in $newC.Keys) {
    $ws.Cells.Item($r, 3).Value = $newC[$r]
}

# --- Delete the leftover VNS column-C values (row 33, rows 242-271) ---
$rowsToDelete = @(33) + (242..271)
foreach ($r in $rowsToDelete) {
    $ws.Cells.Item($r, 3).Clear()
}

# --- Update the visible selection / scroll position --------------------
$ws.Range("C2:C31").Select()

$wb.Save()
